$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Periodo Mora" (arrears period) labels in column E ---
# The four period tags get reshuffled: 1712/1711/1705/1704 -> 1704/1705/1711/1712
$ws.Range("E16").Value = "1704"
$ws.Range("E17").Value = "1705"
$ws.Range("E18").Value = "1711"
$ws.Range("E19").Value = "1712"

# --- Update the matching "Valor Mora" amounts in column F to follow the periods ---
$ws.Range("F16").Value = 29520
$ws.Range("F19").Value = 11808

# --- Reposition the logo image (moved left) ---
$shp = $ws.Shapes.Item("Imagen 2")
$shp.Left = 56.41283464566929
